$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A15").Value = '<jt:escape doublequote="Embedded \"double-quotes\"" backslash="Embedded \\backslash"/>'
